$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style + row height) from the last existing data row (129)
# down into the two new rows being added.
$ws.Range("A129:B129").Copy()
$ws.Range("A130:B131").PasteSpecial(-4122)

# Fill in the new day index / date values
$ws.Cells.Item(130, 1).Value = 129
$ws.Cells.Item(130, 2).Value = 43594

$ws.Cells.Item(131, 1).Value = 130
$ws.Cells.Item(131, 2).Value = 43595

# Match the explicit row height used by the rest of the sheet
$ws.Range("A130:B131").RowHeight = 13.8

# Scroll the view down and select the newly added cells, as in the saved file
$excel.ActiveWindow.ScrollRow = 111
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A130:A131").Select() | Out-Null
